$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new checklist item (row 38) - adds a new shared string and a
# new row referencing it, extending the used range to A1:B38.
$ws.Range("A38").Value = "use correct verb for each request"

# Move the selection/scroll position on to the next empty row, as it was
# left after typing the new entry.
$ws.Range("A39").Select()
